$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-10 (columns A=Região, B=Variável, C=Valor, D=Colocação)
$data = @(
    @("Roraima",  "Diferença 2023/03 - 2022/03",  2.699999999999999,  "1º"),
    @("Amapá",    "Diferença 2023/03 - 2022/03",  1.799999999999999,  "2º"),
    @("Piauí",    "Diferença 2023/03 - 2022/03",  0.7000000000000011, "3º"),
    @("Ceará",    "Diferença 2023/03 - 2022/03",  0.5999999999999996, "4º"),
    @("Amazonas", "Diferença 2023/03 - 2022/03",  0.1999999999999993, "5º"),
    @("Goiás",    "Diferença 2023/03 - 2022/03", -0.1999999999999993, "6º"),
    @("Sergipe",  "Diferença 2023/03 - 2022/03", -2.299999999999999,  "27º"),
    @("Nordeste", "Diferença 2023/03 - 2022/03", -1.199999999999999,  $null),
    @("Brasil",   "Diferença 2023/03 - 2022/03", -0.9999999999999991, $null)
)

# Header row: shift old "Valor" (B1) to C1, insert new B1 and D1 labels
$ws.Range("C1").Value = "Valor"
$ws.Range("B1").Value = "Variável"
$ws.Range("D1").Value = "Colocação"

# Fill rows 2-10
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    if ($row[3] -ne $null) {
        $ws.Cells.Item($r, 4).Value = $row[3]
    } else {
        $ws.Cells.Item($r, 4).Value = ""
    }
}
